$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 112; existing rows 112:131 shift down to 113:132
$ws.Rows.Item(112).Insert()

# Populate the new row 112 with the same constant columns as its neighbours,
# plus the new data values from the diff.
$ws.Range("A112").Value = 6
$ws.Range("B112").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C112").Value = "Metropolitana"
$ws.Range("D112").Value = 45077
$ws.Range("D112").NumberFormat = $ws.Range("D113").NumberFormat
$ws.Range("E112").Value = 13
$ws.Range("F112").Value = 100114007
$ws.Range("G112").Value = "Jengibre"
$ws.Range("H112").Value = "Sin especificar"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 200
$ws.Range("K112").Value = 15000
$ws.Range("L112").Value = 16000
$ws.Range("M112").Value = 15400
$ws.Range("N112").Value = "`$/caja 13 kilos"
$ws.Range("O112").Value = "Perú"
$ws.Range("P112").Value = 1185
$ws.Range("Q112").Value = 13
$ws.Range("R112").Value = "Hortaliza"
